$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the capitalisation of candidate names in column D
$ws.Range("D12").Value = "Ardan Mizanul Khoiri"
$ws.Range("D15").Value = "Mochammad Wafi Nur Jihan"
$ws.Range("D16").Value = "Hoirul Sambudi"
$ws.Range("D17").Value = "Dinda Ayuni"

$ws.Rows.Item(7).AutoFit()

$ws.Range("D10").Select()
